# Add venv row (row 6) to the "C60 db" sheet, and fix a bug where the
# basis value was entered as localized (Farsi) digits, which Excel stores
# as text instead of a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C60 db")

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "برنج ایرانی"
$ws.Cells.Item(6, 3).Value = "۱۳۰"

$ws.Range("D6").Select()
